$wb = $excel.ActiveWorkbook

# Sheet: 展览 (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 112
$ws1.Range("F6").Value = 5351
$ws1.Range("F7").Value = 464
$ws1.Range("F8").Value = 648
$ws1.Range("F9").Value = 923
$ws1.Range("F12").Value = 33
$ws1.Range("F13").Value = 575
$ws1.Range("F14").Value = 22
$ws1.Range("F17").Value = 1787
$ws1.Range("F19").Value = 852
$ws1.Range("F21").Value = 190
$ws1.Range("F22").Value = 316
$ws1.Range("F23").Value = 527
$ws1.Range("F24").Value = 140
$ws1.Range("F28").Value = 2686
$ws1.Range("F29").Value = 175
$ws1.Range("F30").Value = 100
$ws1.Range("F31").Value = 59
$ws1.Range("F32").Value = 98
$ws1.Range("F34").Value = 326
$ws1.Range("F40").Value = 657
$ws1.Range("F42").Value = 49
$ws1.Range("F43").Value = 49

# Sheet: 演出 (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 32

# Sheet: 全部类型 (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 112
$ws4.Range("F6").Value = 32
$ws4.Range("F7").Value = 5351
$ws4.Range("F8").Value = 464
$ws4.Range("F9").Value = 648
$ws4.Range("F12").Value = 923
$ws4.Range("F17").Value = 33
$ws4.Range("F18").Value = 575
$ws4.Range("F19").Value = 22
$ws4.Range("F23").Value = 1787
$ws4.Range("F25").Value = 852
$ws4.Range("F26").Value = 190
$ws4.Range("F27").Value = 316
$ws4.Range("F29").Value = 527
$ws4.Range("F30").Value = 140
$ws4.Range("F32").Value = 2686
$ws4.Range("F33").Value = 175
$ws4.Range("F34").Value = 100
$ws4.Range("F35").Value = 59
$ws4.Range("F36").Value = 98
$ws4.Range("F38").Value = 326
$ws4.Range("F43").Value = 657
$ws4.Range("F45").Value = 49
